$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.956.36"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "3.152.14"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.11"
$ws.Range("E5").Value = "  +1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.60"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.153.55"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.499"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.09"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").Value = "3.671.55"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "64.919.33"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").Value = "3.156.05"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "504.16"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.98"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.16"
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.73"
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.34"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.92"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.53"
$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  +3.66%  "

$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("E34").Value = "  +2.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.48"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.98"
$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0886"
$ws.Range("E37").Value = "  +2.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "476.09"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0414"
$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").Value = "2.991.45"
$ws.Range("E42").Value = "  -3.85%  "

$ws.Range("E43").Value = "  -2.36%  "

$ws.Range("E44").Value = "  -3.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -1.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.25"
$ws.Range("E46").Value = "  -4.21%  "

$ws.Range("D47").Value = "0.0₃0590"
$ws.Range("E47").Value = "  +3.10%  "

$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("E51").Value = "  +14.63%  "
